# Updated cryptos list values (Price / Volume(1h)) per the target diff.
# Price cells that look like plain numbers are entered with a leading
# apostrophe so Excel stores them as text (matching the source format,
# e.g. "1.003", "41.87") instead of re-parsing them as numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.430.24'
$ws.Range("E2").Value = '  +1.69%  '
$ws.Range("D3").Value = '1.826.89'
$ws.Range("E3").Value = '  +1.81%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''316.71'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '''0.5158'
$ws.Range("E7").Value = '  -3.17%  '
$ws.Range("D8").Value = '''0.3907'
$ws.Range("E8").Value = '  -1.36%  '
$ws.Range("D9").Value = '''0.07640'
$ws.Range("E9").Value = '  +2.38%  '
$ws.Range("D10").Value = '''41.87'
$ws.Range("E10").Value = '  +1.13%  '
$ws.Range("D11").Value = '''1.108'
$ws.Range("E11").Value = '  +2.18%  '
$ws.Range("D12").Value = '''21.04'
$ws.Range("E12").Value = '  +3.31%  '
$ws.Range("D13").Value = '''6.279'
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").Value = '''7.539'
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '1.822.76'
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("D17").Value = '''93.33'
$ws.Range("E17").Value = '  +5.57%  '
$ws.Range("D18").Value = '''0.00001080'
$ws.Range("E18").Value = '  +1.94%  '
$ws.Range("D19").Value = '''0.06668'
$ws.Range("E19").Value = '  +1.46%  '
$ws.Range("D20").Value = '''17.66'
$ws.Range("E20").Value = '  +2.59%  '
$ws.Range("D21").Value = '''1.001'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").Value = '''6.185'
$ws.Range("D23").Value = '28.475.98'
$ws.Range("E23").Value = '  +1.76%  '
$ws.Range("D24").Value = '''11.12'
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("D25").Value = '''2.255'
$ws.Range("E25").Value = '  +7.99%  '
$ws.Range("D26").Value = '''156.90'
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").Value = '''20.61'
$ws.Range("E27").Value = '  +2.02%  '
$ws.Range("D28").Value = '2.034.81'
$ws.Range("E28").Value = '  +1.95%  '
$ws.Range("D29").Value = '''2.395'
$ws.Range("E29").Value = '  +4.03%  '
$ws.Range("D30").Value = '''124.79'
$ws.Range("E30").Value = '  +2.23%  '
$ws.Range("D31").Value = '''1.119'
$ws.Range("E31").Value = '  +2.39%  '
$ws.Range("D32").Value = '''0.1084'
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("D33").Value = '''5.659'
$ws.Range("E33").Value = '  +2.76%  '
$ws.Range("D34").Value = '''3.664'
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D35").Value = '''0.06995'
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("D36").Value = '''0.2224'
$ws.Range("E36").Value = '  +0.33%  '
$ws.Range("D37").Value = '''8.974'
$ws.Range("E37").Value = '  +6.86%  '
$ws.Range("D38").Value = '''0.02322'
$ws.Range("E38").Value = '  +2.14%  '
$ws.Range("D39").Value = '''5.130'
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("D40").Value = '''0.6280'
$ws.Range("E40").Value = '  +2.56%  '
$ws.Range("D41").Value = '''11.21'
$ws.Range("E41").Value = '  -0.33%  '
$ws.Range("D42").Value = '''1.184'
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '''1.397'
$ws.Range("E44").Value = '  -1.27%  '
$ws.Range("D45").Value = '''13.39'
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("D46").Value = '''0.5896'
$ws.Range("E46").Value = '  +3.08%  '
$ws.Range("D47").Value = '''3.712'
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("D48").Value = '''124.15'
$ws.Range("E48").Value = '  -0.79%  '
$ws.Range("D49").Value = '''1.976'
$ws.Range("E49").Value = '  +2.90%  '
$ws.Range("D50").Value = '''1.199'
$ws.Range("E50").Value = '  +1.52%  '
$ws.Range("D51").Value = '''0.06929'
$ws.Range("E51").Value = '  +1.83%  '
